# MCHP Football Fund — record August ,18 (column G) payments that just
# came in from Aniket (row 5), Pravin (row 18) and Sadik (row 22).
# The dependent SUM/balance formulas elsewhere on the sheet recalculate
# automatically from these three inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 500
$ws.Range("G18").Value = 500
$ws.Range("G22").Value = 500

# Leave the view scrolled down to where the edits were made, with G30
# as the active cell (best-effort — scroll position is host UI state).
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("G30").Select()
